# Insert a new "Development Tools" slide as the 3rd slide (after "Project
# Objective", before "Linear Regression"), using the same Title+Content
# layout used by the other content slides.

$p = $ppt.ActivePresentation

# ppLayoutText = 2 -> resolves to the "Title and Content" custom layout,
# same layout already used by every other content slide in this deck.
$s = $p.Slides.Add(3, 2)

# Title placeholder.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Development Tools"

# Body / content placeholder - build it run-by-run so the paragraph/run
# breaks match how PowerPoint splits text around the spell-checked words
# ("Javascript", "sklearn", "numpy").
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Front end: HTML, CSS, "
[void]$body.InsertAfter("Javascript")
[void]$body.InsertAfter(" (Node)`r")
[void]$body.InsertAfter("Backend: Python, flask, ")
[void]$body.InsertAfter("sklearn")
[void]$body.InsertAfter(", ")
[void]$body.InsertAfter("numpy")
